# Auto-generated edit script: update cryptos price list values
# (refreshed data pulled on Tue Aug  1 20:47:18 UTC 2023 by GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.187.64"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
$ws.Range("D3").Value = "1.849.26"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "'245.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "

# Row 6
$ws.Range("D6").Value = "'0.6979"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.14%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.07713"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("D9").Value = "'0.3062"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "

# Row 10
$ws.Range("D10").Value = "'23.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "

# Row 11
$ws.Range("D11").Value = "'0.07821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("D12").Value = "'93.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.16%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.846.57"
$ws.Range("E13").Value = "  -1.06%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.125"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

# Row 15
$ws.Range("D15").Value = "'0.6850"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("D16").Value = "'6.635"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.172.97"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000008297"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

# Row 19
$ws.Range("D19").Value = "'241.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "

# Row 20
$ws.Range("D20").Value = "2.087.27"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").Value = "'12.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").Value = "'7.515"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").Value = "'0.1509"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("D26").Value = "'159.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "

# Row 27
$ws.Range("D27").Value = "'8.826"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "

# Row 28
$ws.Range("E28").Value = "  -1.29%  "

# Row 29
$ws.Range("D29").Value = "'1.547"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "

# Row 30
$ws.Range("D30").Value = "'4.227"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "

# Row 31
$ws.Range("E31").Value = "  -0.82%  "

# Row 32
$ws.Range("D32").Value = "'1.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("D33").Value = "'0.05120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "

# Row 34
$ws.Range("D34").Value = "'0.7946"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "

# Row 35
$ws.Range("D35").Value = "'1.866"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "

# Row 36
$ws.Range("D36").Value = "'1.147"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "

# Row 37
$ws.Range("D37").Value = "'2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
$ws.Range("D38").Value = "1.319.49"
$ws.Range("E38").Value = "  +7.70%  "

# Row 39
$ws.Range("D39").Value = "'0.01873"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.87%  "

# Row 40
$ws.Range("D40").Value = "'2.713"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.30%  "

# Row 41
$ws.Range("D41").Value = "'0.9480"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.60%  "

# Row 42
$ws.Range("D42").Value = "'6.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.45%  "

# Row 43
$ws.Range("D43").Value = "'107.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "

# Row 44
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "

# Row 45
$ws.Range("D45").Value = "'9.713"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.13%  "

# Row 46
$ws.Range("D46").Value = "1.988.26"
$ws.Range("E46").Value = "  -0.94%  "

# Row 47
$ws.Range("D47").Value = "'0.5182"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").Value = "'64.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.42%  "

# Row 49
$ws.Range("D49").Value = "'1.763"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "

# Row 50
$ws.Range("D50").Value = "'0.00000000118"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.88%  "

# Row 51
$ws.Range("D51").Value = "'6.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "

